$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (account 005142592 / ALBERTO / 345310)
$ws.Rows.Item(2).Delete()

# Insert a new row for CASSIO right after the row for ANA (004267119 / 186.82),
# i.e. just before the GUSTAVO (005591536 / 129.8) row, which is now row 28.
$ws.Rows.Item(28).Insert()

# Account number must stay text (keep the leading zero), use an apostrophe
# prefix so Excel stores it as a plain text value, matching the inlineStr cells
# used throughout the rest of the sheet. ClearFormats afterwards so the cell
# doesn't pick up a stray "quote prefix" number format and instead matches the
# default (unstyled) look of the surrounding data cells.
$ws.Range("A28").Value = "'004508526"
$ws.Range("A28").ClearFormats()
$ws.Range("B28").Value = "CASSIO"
$ws.Range("C28").Value = 153.7
